$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, shifting existing rows 198:235 down to 199:236.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly price record.
$ws.Cells.Item(198, 1).Value = 10
$ws.Cells.Item(198, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(198, 3).Value = "La Araucanía"
$ws.Cells.Item(198, 4).Value = 45218
$ws.Cells.Item(198, 5).Value = 9
$ws.Cells.Item(198, 6).Value = "Fruta"
$ws.Cells.Item(198, 7).Value = 100107
$ws.Cells.Item(198, 8).Value = "Otros"
$ws.Cells.Item(198, 9).Value = 100107002
$ws.Cells.Item(198, 10).Value = "Chirimoya"
$ws.Cells.Item(198, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(198, 12).Value = "Primera"
$ws.Cells.Item(198, 13).Value = 260
$ws.Cells.Item(198, 14).Value = 2500
$ws.Cells.Item(198, 15).Value = 2600
$ws.Cells.Item(198, 16).Value = 2569
$ws.Cells.Item(198, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(198, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(198, 19).Value = 2569
$ws.Cells.Item(198, 20).Value = 1
